$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Week15"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'Lundo’s Legends'
$ws.Range("C2").Value = 131
$ws.Range("D2").Value = 57
$ws.Range("E2").Value = 74
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = -2
$ws.Range("H2").Value = 2.5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'EL Onće'
$ws.Range("C3").Value = 122
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = 57
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = -1

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'Epic7'
$ws.Range("C4").Value = 121
$ws.Range("D4").Value = 51.5
$ws.Range("E4").Value = 69.5
$ws.Range("F4").Value = -3.5
$ws.Range("G4").Value = -2
$ws.Range("H4").Value = -1.5

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Samsquanches'
$ws.Range("C5").Value = 107
$ws.Range("D5").Value = 66
$ws.Range("E5").Value = 41
$ws.Range("F5").Value = 6.5
$ws.Range("G5").Value = 10.5
$ws.Range("H5").Value = -4

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'Splitfinger Skadoosh'
$ws.Range("C6").Value = 89.5
$ws.Range("D6").Value = 45
$ws.Range("E6").Value = 44.5
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = -2.5
$ws.Range("H6").Value = 9.5

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'confusion'
$ws.Range("C7").Value = 88
$ws.Range("D7").Value = 47.5
$ws.Range("E7").Value = 40.5
$ws.Range("F7").Value = 1.5
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.5

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 'SmokeWalkers'
$ws.Range("C8").Value = 83
$ws.Range("D8").Value = 49
$ws.Range("E8").Value = 34
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = -1
$ws.Range("H8").Value = 1.5

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'GOD WILLS IT'
$ws.Range("C9").Value = 83
$ws.Range("D9").Value = 43.5
$ws.Range("E9").Value = 39.5
$ws.Range("F9").Value = -6
$ws.Range("G9").Value = -2.5
$ws.Range("H9").Value = -3.5

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'DJ''s Quality Team'
$ws.Range("C10").Value = 83
$ws.Range("D10").Value = 44
$ws.Range("E10").Value = 39
$ws.Range("F10").Value = -6
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = -8

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'rainmaker'
$ws.Range("C11").Value = 79
$ws.Range("D11").Value = 20.5
$ws.Range("E11").Value = 58.5
$ws.Range("F11").Value = -1
$ws.Range("G11").Value = -3
$ws.Range("H11").Value = 2

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'Swampnuts'
$ws.Range("C12").Value = 77.5
$ws.Range("D12").Value = 44
$ws.Range("E12").Value = 33.5
$ws.Range("F12").Value = -13
$ws.Range("G12").Value = -9.5
$ws.Range("H12").Value = -3.5

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'MillerTime'
$ws.Range("C13").Value = 72
$ws.Range("D13").Value = 32
$ws.Range("E13").Value = 40
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = -0.5
$ws.Range("H13").Value = 9.5

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'Corbin Copy'
$ws.Range("C14").Value = 64
$ws.Range("D14").Value = 41
$ws.Range("E14").Value = 23
$ws.Range("F14").Value = -6.5
$ws.Range("G14").Value = -6
$ws.Range("H14").Value = -0.5

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'Mac'
$ws.Range("C15").Value = 60
$ws.Range("D15").Value = 24
$ws.Range("E15").Value = 36
$ws.Range("F15").Value = -2
$ws.Range("G15").Value = 1.5
$ws.Range("H15").Value = -3.5

